# "Generate Report for Handoff" -- adds two new localization-status rows
# (image dependency handoffs) and promotes the existing row from a plain
# ".md" handoff to a full zh-cn/de-de xlf handoff, across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$commitE2e   = "0382163a91ac6e9f262b70c2ce9842bf657a0af7"
$commitZhCn  = "e873cba4933fdbff3181360e4842855e4ab4b978"
$commitDeDe  = "6469762daeca757c878185ae03f73b54649792d8"
$xlfHash     = "6b9ed68ec09722ab816be9693e192e9339406eaf"

$mdFile   = "fa47e20e-7d45-4e4f-8606-15b186488cec.md"
$png1     = "25bd4e33-8b9b-423d-8cd9-dad4ac932d88.png"
$png2     = "7c492253-f256-466e-acb6-df8c7299bf0c.png"
$png1Tgt  = "586ad7dcd0e5287f1c7bb7e1bcd0fe9d6ddd2ee2.png"
$png2Tgt  = "1e846f58beb49082783cd9937bec80ce60d5a3eb.png"
$xlfZh    = "$mdFile.$xlfHash.zh-cn.xlf"
$xlfDe    = "$mdFile.$xlfHash.de-de.xlf"

$readyStatus = "Ready for handoff"
$dateOverview = "2016-50-20 12:50:19"
$dateHandoffZh = "2016-03-20 12:50:16"
$dateHandoffDe = "2016-03-20 12:50:19"
$zeroDate = "0001-01-01 00:00:00"
$depFrom = "e2e\$mdFile"

function Set-Hyperlink($ws, $cellAddr, $address, $text) {
    $ws.Range($cellAddr).Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $address, "", "", $text) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Existing row 2 becomes the first png dependency entry.
$ov.Range("B2").Value = $readyStatus
$ov.Range("C2").Value = $readyStatus
$ov.Range("D2").Value = $dateOverview
Set-Hyperlink $ov "A2" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$png1" $png1

# New row 3: second png dependency entry.
$ov.Range("B3").Value = $readyStatus
$ov.Range("C3").Value = $readyStatus
$ov.Range("D3").Value = $dateOverview
Set-Hyperlink $ov "A3" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$png2" $png2

# New row 4: the markdown file itself.
$ov.Range("B4").Value = $readyStatus
$ov.Range("C4").Value = $readyStatus
$ov.Range("D4").Value = $dateOverview
Set-Hyperlink $ov "A4" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$mdFile" $mdFile

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = ".png"
$zh.Range("E2").Value = $dateHandoffZh
$zh.Range("H2").Value = $zeroDate
$zh.Range("I2").Value = "IsDependency"
$zh.Range("J2").Value = $depFrom
Set-Hyperlink $zh "A2" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$png1" $png1
Set-Hyperlink $zh "B2" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$png1" ".png"
Set-Hyperlink $zh "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitZhCn/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png1Tgt" $png1Tgt

$zh.Range("A3").Value = $png2
$zh.Range("B3").Value = ".png"
$zh.Range("C3").Value = $readyStatus
$zh.Range("D3").Value = $png2Tgt
$zh.Range("E3").Value = $dateHandoffZh
$zh.Range("H3").Value = $zeroDate
$zh.Range("I3").Value = "IsDependency"
$zh.Range("J3").Value = $depFrom
$zh.Range("A3").Style = "Normal"
$zh.Range("B3").Style = "Normal"
$zh.Range("D3").Style = "Normal"
Set-Hyperlink $zh "A3" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$png2" $png2
Set-Hyperlink $zh "B3" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$png2" ".png"
Set-Hyperlink $zh "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitZhCn/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png2Tgt" $png2Tgt

$zh.Range("A4").Value = $mdFile
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = $readyStatus
$zh.Range("D4").Value = $xlfZh
$zh.Range("E4").Value = $dateHandoffZh
$zh.Range("H4").Value = $zeroDate
$zh.Range("I4").Value = "Include"
$zh.Range("A4").Style = "Normal"
$zh.Range("B4").Style = "Normal"
$zh.Range("D4").Style = "Normal"
Set-Hyperlink $zh "A4" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$mdFile" $mdFile
Set-Hyperlink $zh "B4" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$mdFile" ".md"
Set-Hyperlink $zh "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitZhCn/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZh" $xlfZh

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = ".png"
$de.Range("E2").Value = $dateHandoffDe
$de.Range("H2").Value = $zeroDate
$de.Range("I2").Value = "IsDependency"
$de.Range("J2").Value = $depFrom
Set-Hyperlink $de "A2" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$png1" $png1
Set-Hyperlink $de "B2" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$png1" ".png"
Set-Hyperlink $de "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitDeDe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png1Tgt" $png1Tgt

$de.Range("A3").Value = $png2
$de.Range("B3").Value = ".png"
$de.Range("C3").Value = $readyStatus
$de.Range("D3").Value = $png2Tgt
$de.Range("E3").Value = $dateHandoffDe
$de.Range("H3").Value = $zeroDate
$de.Range("I3").Value = "IsDependency"
$de.Range("J3").Value = $depFrom
$de.Range("A3").Style = "Normal"
$de.Range("B3").Style = "Normal"
$de.Range("D3").Style = "Normal"
Set-Hyperlink $de "A3" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$png2" $png2
Set-Hyperlink $de "B3" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$png2" ".png"
Set-Hyperlink $de "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitDeDe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png2Tgt" $png2Tgt

$de.Range("A4").Value = $mdFile
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = $readyStatus
$de.Range("D4").Value = $xlfDe
$de.Range("E4").Value = $dateHandoffDe
$de.Range("H4").Value = $zeroDate
$de.Range("I4").Value = "Include"
$de.Range("A4").Style = "Normal"
$de.Range("B4").Style = "Normal"
$de.Range("D4").Style = "Normal"
Set-Hyperlink $de "A4" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$mdFile" $mdFile
Set-Hyperlink $de "B4" "https://github.com/OpenLocalizationTest/oltest/blob/$commitE2e/e2e/$mdFile" ".md"
Set-Hyperlink $de "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitDeDe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDe" $xlfDe

Write-Output "done"
